$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8081873059272766
$ws.Range("B1").Value = 1.518516302108765
$ws.Range("C1").Value = 5.935915946960449
$ws.Range("D1").Value = 3.115354776382446
$ws.Range("E1").Value = 1.46107006072998
